$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D21:F21").ClearContents()
$ws.Range("F22").ClearContents()

$ws.Range("G27").Select()
